$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header labels: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $colLeft = $i + 1        # columns A..J
    $colRight = $i + 12      # columns L..U
    $ws.Cells.Item(1, $colLeft).Value = $baseNames[$i] + "_FV2310"
    $ws.Cells.Item(1, $colRight).Value = $baseNames[$i] + "_FV2404"
}

# --- Turn the used range into an Excel Table ("Table1") ---
$dataRange = $ws.Range("A1:U76")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- Freeze the header row (split after row 1, keep column A as top-left of scrolling pane) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
